$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns: "Ma dinh danh vi tri" (AA) and "Tinh trang hieu luc vi tri" (AB) ---
$ws.Range("AA3").Value = "Mã định danh vị trí"
$ws.Range("AB3").Value = "Tình trạng hiệu lực vị trí"

# Merge header cells vertically (row3:row4), matching the other header columns
$ws.Range("AA3:AA4").Merge()
$ws.Range("AB3:AB4").Merge()

# Style for AA3:AA4 -> bold font, thin left border only, centered, wrap text
$rAA = $ws.Range("AA3:AA4")
$rAA.Font.Bold = $true
$rAA.HorizontalAlignment = -4108
$rAA.WrapText = $true
$rAA.Borders.Item(7).LineStyle = 1

# Style for AB3:AB4 -> bold font, no border, centered, wrap text
$rAB = $ws.Range("AB3:AB4")
$rAB.Font.Bold = $true
$rAB.HorizontalAlignment = -4108
$rAB.WrapText = $true

# --- View state tweaks (best-effort; match new scroll/selection/zoom) ---
$win = $excel.ActiveWindow
$win.Zoom = 100
$win.ScrollRow = 1
$win.ScrollColumn = 12
$ws.Range("Z11").Select()
